$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two stray local-file RSS feed rows (payloadspace / satnews) ---
# Only columns A (feed URL) and B (category) held data on these rows; the
# keyword/category lookup table in columns C/D on the same rows is untouched.
$ws.Range("A38:B38").ClearContents()
$ws.Range("A39:B39").ClearContents()

# --- Append four new keyword -> category rows to the lookup table (C/D) ---
$ws.Range("C167").Value = "Hale"
$ws.Range("D167").Value = "Competitors"

$ws.Range("C168").Value = "Earth observation"
$ws.Range("D168").Value = "Satellite Operations"

$ws.Range("C169").Value = "Data"
$ws.Range("D169").Value = "Satellite Operations"

$ws.Range("C170").Value = "Weather"
$ws.Range("D170").Value = "Satellite Operations"

# Keep the new rows' height consistent with the rest of the lookup table
$ws.Range("C167:D170").RowHeight = 15.75

# --- Match the new view state recorded in the workbook ---
$ws.Activate()
$ws.Range("B2").Select()
